{"js": "// Replace each two-digit-division expression with its updated value.\n// Mapping preserves document order and is derived from the commit diff;\n// since every \"old\" string is unique in the document and no replacement\n// value collides with a not-yet-processed \"old\" value, sequential\n// search-and-replace is safe and deterministic.\nconst replacements = [[\"30\u00f76=\", \"17\u00f77=\"], [\"10\u00f76=\", \"11\u00f74=\"], [\"82\u00f74=\", \"42\u00f72=\"], [\"80\u00f76=\", \"64\u00f79=\"], [\"31\u00f77=\", \"43\u00f75=\"], [\"24\u00f76=\", \"28\u00f76=\"], [\"79\u00f72=\", \"14\u00f75=\"], [\"73\u00f75=\", \"99\u00f79=\"], [\"77\u00f74=\", \"29\u00f73=\"], [\"80\u00f75=\", \"63\u00f73=\"], [\"11\u00f72=\", \"24\u00f75=\"], [\"61\u00f75=\", \"90\u00f75=\"], [\"62\u00f76=\", \"56\u00f72=\"], [\"56\u00f79=\", \"51\u00f79=\"], [\"30\u00f75=\", \"60\u00f74=\"], [\"91\u00f75=\", \"44\u00f76=\"], [\"39\u00f72=\", \"94\u00f72=\"], [\"10\u00f79=\", \"55\u00f72=\"], [\"55\u00f74=\", \"63\u00f75=\"], [\"75\u00f74=\", \"51\u00f79=\"], [\"41\u00f75=\", \"92\u00f74=\"], [\"12\u00f78=\", \"55\u00f79=\"], [\"44\u00f73=\", \"22\u00f77=\"], [\"76\u00f77=\", \"45\u00f78=\"], [\"35\u00f77=\", \"21\u00f73=\"]];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text \"${oldText}\" to replace.`);\n  }\n\n  // Replace the first (and, per the source document, only) occurrence.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update each two-digit division expression to its new value.\n# The mapping mirrors the commit diff; every \"old\" string is unique in\n# the document, so a literal, non-wildcard Find/Replace scoped to\n# \"replace one\" occurrence keeps each substitution precise and avoids\n# any unintended cascading replacements.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"30\u00f76=\", \"17\u00f77=\"),\n    @(\"10\u00f76=\", \"11\u00f74=\"),\n    @(\"82\u00f74=\", \"42\u00f72=\"),\n    @(\"80\u00f76=\", \"64\u00f79=\"),\n    @(\"31\u00f77=\", \"43\u00f75=\"),\n    @(\"24\u00f76=\", \"28\u00f76=\"),\n    @(\"79\u00f72=\", \"14\u00f75=\"),\n    @(\"73\u00f75=\", \"99\u00f79=\"),\n    @(\"77\u00f74=\", \"29\u00f73=\"),\n    @(\"80\u00f75=\", \"63\u00f73=\"),\n    @(\"11\u00f72=\", \"24\u00f75=\"),\n    @(\"61\u00f75=\", \"90\u00f75=\"),\n    @(\"62\u00f76=\", \"56\u00f72=\"),\n    @(\"56\u00f79=\", \"51\u00f79=\"),\n    @(\"30\u00f75=\", \"60\u00f74=\"),\n    @(\"91\u00f75=\", \"44\u00f76=\"),\n    @(\"39\u00f72=\", \"94\u00f72=\"),\n    @(\"10\u00f79=\", \"55\u00f72=\"),\n    @(\"55\u00f74=\", \"63\u00f75=\"),\n    @(\"75\u00f74=\", \"51\u00f79=\"),\n    @(\"41\u00f75=\", \"92\u00f74=\"),\n    @(\"12\u00f78=\", \"55\u00f79=\"),\n    @(\"44\u00f73=\", \"22\u00f77=\"),\n    @(\"76\u00f77=\", \"45\u00f78=\"),\n    @(\"35\u00f77=\", \"21\u00f73=\"),\n)\n\n$wdReplaceOne = 1\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, $wdReplaceOne)\n}\n\n"}
